$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.454.80"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.591.71"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.72"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.47"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.814.13"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.583.42"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.01"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.08"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.97%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.62%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "206.56"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.03"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.63%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.91%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.84"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.02"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.18%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0503"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.75%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.00%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.280.22"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.50"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.08%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +12.48%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.74%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.41%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.66%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.39"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.14"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.768"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "62.14"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.726.43"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "88.65"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.97%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.30%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0512"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.399"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.44%  "
